$d = $word.ActiveDocument
Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
$d.Content.Find.Execute("Google ", $true, $false, $false, $false, $false, $true, 1, $false, "TESTING", 2)
